# atualização da rodada 7
# Rodada 7 da premier league
#
# Adds the results of round 7 to the "cartoes" (yellow cards) table.
# For each club (rows 2-21) the counters below are updated:
#   B partidas jogadas              (matches played)
#   C partidas como mandante        (matches played at home)
#   D cartoes amarelos mand         (yellow cards received at home)
#   E média de cartões amar. mand   (= D / C)
#   F partidas como visitante       (matches played away)
#   G cartoes amarelos visit        (yellow cards received away)
#   H média de cart amar.visit      (= G / F)
#   I total de cartoes amar.        (= D + G)
#   J Média cart. amar              (= I / B)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New totals (after round 7) for B, C, D, F, G per row; E, H, I, J are derived.
$data = @{
    2  = @(7, 4, 7,  3, 10)
    3  = @(7, 4, 6,  3, 8)
    4  = @(7, 3, 10, 4, 9)
    5  = @(7, 4, 7,  3, 5)
    6  = @(7, 4, 10, 3, 5)
    7  = @(7, 4, 14, 3, 13)
    8  = @(7, 4, 10, 3, 7)
    9  = @(7, 4, 8,  3, 7)
    10 = @(7, 4, 6,  3, 6)
    11 = @(7, 3, 10, 4, 11)
    12 = @(7, 3, 9,  4, 11)
    13 = @(7, 4, 10, 3, 9)
    14 = @(7, 3, 7,  4, 7)
    15 = @(7, 3, 11, 4, 13)
    16 = @(7, 3, 10, 4, 7)
    17 = @(7, 3, 6,  4, 13)
    18 = @(7, 3, 8,  4, 13)
    19 = @(7, 3, 8,  4, 10)
    20 = @(7, 4, 9,  3, 7)
    21 = @(7, 3, 7,  4, 15)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $b = $vals[0]
    $c = $vals[1]
    $d = $vals[2]
    $f = $vals[3]
    $g = $vals[4]

    $e = $d / $c
    $h = $g / $f
    $i = $d + $g
    $j = $i / $b

    $ws.Cells.Item($row, 2).Value  = $b   # B
    $ws.Cells.Item($row, 3).Value  = $c   # C
    $ws.Cells.Item($row, 4).Value  = $d   # D
    $ws.Cells.Item($row, 5).Value  = $e   # E
    $ws.Cells.Item($row, 6).Value  = $f   # F
    $ws.Cells.Item($row, 7).Value  = $g   # G
    $ws.Cells.Item($row, 8).Value  = $h   # H
    $ws.Cells.Item($row, 9).Value  = $i   # I
    $ws.Cells.Item($row, 10).Value = $j   # J
}
